# point & coupon page added, database updated (please add coupon & point)
#
# Adds to Sheet1:
#   - A small manager/seller login table in columns L:N (rows 4-6)
#   - A new "point" table block (rows 75-81)
#   - A new "coupon" table block (rows 85-92), with two free-text note rows
#     (83 and 91-92)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- manager / seller credential table (columns L:N, rows 4-6) ---
$ws.Range("M4").Value = "id"
$ws.Range("N4").Value = "passwd"

$ws.Range("L5").Value = "관리자"
$ws.Range("M5").Value = "manager"
$ws.Range("N5").Value = 1234

$ws.Range("L6").Value = "판매자"
$ws.Range("M6").Value = "seller"
$ws.Range("N6").Value = 1234

# --- "point" table block ---
# Table-name cells (column C) reuse the formatting already used by the
# other table-name cells (e.g. C3 "userinfo").
$ws.Range("C3").Copy()
$ws.Range("C75").PasteSpecial(-4122)
$ws.Range("C83").PasteSpecial(-4122)

# Field-name cells (column D) reuse the formatting already used by the
# other field-name cells (e.g. D5 "id").
$ws.Range("D5").Copy()
$ws.Range("D77:D80").PasteSpecial(-4122)

$ws.Range("C75").Value = "point"
$ws.Range("E76").Value = "CREATE TABLE ``point`` ("

$ws.Range("D77").Value = "아이디"
$ws.Range("F77").Value = "``id`` VARCHAR(50) NULL DEFAULT NULL,"

$ws.Range("D78").Value = "주문번호"
$ws.Range("F78").Value = "``order_num`` INT(11) NULL DEFAULT NULL,"

$ws.Range("D79").Value = "포인트금액"
$ws.Range("F79").Value = "``point`` VARCHAR(50) NULL DEFAULT NULL,"

$ws.Range("D80").Value = "날짜"
$ws.Range("F80").Value = "``date`` VARCHAR(50) NULL DEFAULT NULL"

$ws.Range("E81").Value = ")"

$ws.Range("C83").Value = "주문하는 페이지에서 (insert)입력, 마이페이지에서 id별로 (select)보여주기, 주문취소시 삭제 - 주문페이지 제작시 주분번호를 생성 필요함,"

# --- "coupon" table block ---
$ws.Range("C3").Copy()
$ws.Range("C85").PasteSpecial(-4122)
$ws.Range("C91").PasteSpecial(-4122)
$ws.Range("C92").PasteSpecial(-4122)

$ws.Range("D5").Copy()
$ws.Range("D86:D89").PasteSpecial(-4122)

$ws.Range("C85").Value = "coupon"
$ws.Range("E85").Value = "CREATE TABLE ``coupon`` ("

$ws.Range("D86").Value = "아이디"
$ws.Range("F86").Value = "``id`` VARCHAR(50) NULL DEFAULT NULL,"

$ws.Range("D87").Value = "쿠폰종류"
$ws.Range("F87").Value = "``coupontype`` INT(11) NULL DEFAULT NULL,"

$ws.Range("D88").Value = "유효일자"
$ws.Range("F88").Value = "``end_date`` INT(11) NULL DEFAULT NULL,"

$ws.Range("D89").Value = "생성일자"
$ws.Range("F89").Value = "``start_date`` INT(11) NULL DEFAULT NULL"

$ws.Range("E90").Value = ")"

$ws.Range("C91").Value = "쿠폰 종류를 번호에 따라 정해두고, 판매자 페이지에서 발급하는 기능 추가하기 ( 판매자 페이지에서 coupon table에 insert)"
$ws.Range("C92").Value = "쿠폰 종류: 1 (5% 할인) 2 (10% 할인) 3 (30% 할인) 4 (1000원 할인)"

# --- move the selection to where the new data now lives ---
$ws.Range("N5").Select()
